$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.091.80"
$ws.Range("E2").Value = "  +0.88%  "
$ws.Range("D3").Value = "'2.642.57"
$ws.Range("E3").Value = "  +4.16%  "
$ws.Range("D5").Value = "'518.04"
$ws.Range("E5").Value = "  +2.05%  "
$ws.Range("D6").Value = "'144.53"
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("D7").Value = "'0.996"
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("D8").Value = "'0.570"
$ws.Range("E8").Value = "  +1.12%  "
$ws.Range("D9").Value = "'2.669.36"
$ws.Range("E9").Value = "  +5.08%  "
$ws.Range("D10").Value = "'6.27"
$ws.Range("E10").Value = "  +1.41%  "
$ws.Range("E11").Value = "  +3.46%  "
$ws.Range("E12").Value = "  +2.05%  "
$ws.Range("E13").Value = "  -1.53%  "
$ws.Range("D14").Value = "'3.114.21"
$ws.Range("E14").Value = "  +4.52%  "
$ws.Range("D15").Value = "'59.022.26"
$ws.Range("E15").Value = "  +0.79%  "
$ws.Range("D16").Value = "'21.01"
$ws.Range("E16").Value = "  +1.57%  "
$ws.Range("D17").Value = "'0.0000138"
$ws.Range("E17").Value = "  +1.81%  "
$ws.Range("D18").Value = "'2.668.46"
$ws.Range("E18").Value = "  +5.17%  "
$ws.Range("D19").Value = "'350.35"
$ws.Range("E19").Value = "  +4.37%  "
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("D21").Value = "'10.38"
$ws.Range("E21").Value = "  +3.11%  "
$ws.Range("E22").Value = "  +4.21%  "
$ws.Range("D23").Value = "'0.997"
$ws.Range("E23").Value = "  -0.32%  "
$ws.Range("D24").Value = "'62.09"
$ws.Range("E24").Value = "  +3.24%  "
$ws.Range("E25").Value = "  +3.08%  "
$ws.Range("D26").Value = "'0.993"
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("E27").Value = "  +0.93%  "
$ws.Range("D28").Value = "'0.0₃0807"
$ws.Range("E28").Value = "  +2.41%  "
$ws.Range("D29").Value = "'7.15"
$ws.Range("E29").Value = "  +2.82%  "
$ws.Range("D30").Value = "'0.997"
$ws.Range("E30").Value = "  -0.26%  "
$ws.Range("E31").Value = "  +8.52%  "
$ws.Range("E32").Value = "  +2.97%  "
$ws.Range("E33").Value = "  +3.06%  "
$ws.Range("D34").Value = "'150.26"
$ws.Range("E34").Value = "  +0.40%  "
$ws.Range("D35").Value = "'0.973"
$ws.Range("E35").Value = "  +4.29%  "
$ws.Range("E36").Value = "  +2.98%  "
$ws.Range("E37").Value = "  +2.62%  "
$ws.Range("D38").Value = "'36.67"
$ws.Range("E38").Value = "  +1.70%  "
$ws.Range("D39").Value = "'0.843"
$ws.Range("E39").Value = "  +2.20%  "
$ws.Range("D40").Value = "'3.72"
$ws.Range("E40").Value = "  +5.75%  "
$ws.Range("E41").Value = "  +0.84%  "
$ws.Range("D42").Value = "'279.35"
$ws.Range("E42").Value = "  -1.50%  "
$ws.Range("D43").Value = "'0.614"
$ws.Range("E43").Value = "  +2.01%  "
$ws.Range("E44").Value = "  -0.68%  "
$ws.Range("D45").Value = "'0.993"
$ws.Range("D46").Value = "'19.67"
$ws.Range("E46").Value = "  +5.81%  "
$ws.Range("E47").Value = "  -0.49%  "
$ws.Range("E48").Value = "  +2.11%  "
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("D50").Value = "'1.992.94"
$ws.Range("E50").Value = "  +5.27%  "
$ws.Range("E51").Value = "  +3.30%  "
